$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2Emily"
$ws.Range("B2").Value = "sid021111"
$ws.Range("C2").Value = "spw021111"
$ws.Range("D2").Value = "A Mobile App to search for the food addictives by computer vision"

$ws.Range("A3").Value = "2Fanny"
$ws.Range("B3").Value = "sid022222"
$ws.Range("C3").Value = "spw022222"
$ws.Range("D3").Value = "Python Online Practice System"

$ws.Range("E1").Value = "credit"
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 6

$ws.Range("E3").Select() | Out-Null
